$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-57 down to 8-58
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new data record
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44503
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100102
$ws.Range("H7").Value = "Cítricos"
$ws.Range("I7").Value = 100102006
$ws.Range("J7").Value = "Pomelo"
$ws.Range("K7").Value = "Start Ruby"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 580
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7500
$ws.Range("P7").Value = 7259
$ws.Range("Q7").Value = "`$/caja 14 kilos granel"
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 518
$ws.Range("T7").Value = 14

Write-Host "Done"
